# Segunda edición: agrega la columna A (filas 9-12) con nuevas etiquetas
# y actualiza la celda seleccionada en la hoja.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = "ADAD"
$ws.Range("A10").Value = "Ojalá lo note"
$ws.Range("A11").Value = "El gitHub"
$ws.Range("A12").Value = "xdd"

$ws.Range("B11").Select()
